{"js": "// Change: \"Below are the cross-validation results for all the three\n// algorithms.\" -> \"Below are the cross-validation results for both the\n// algorithms.\" (the sentence that immediately precedes the algorithm\n// comparison table). The document contains another, unrelated sentence\n// (\"Columns mentioned in all the three categories ...\") that must stay\n// untouched, so we first locate the whole target sentence (which is\n// unique) and then scope the replacement search to that sentence only.\n\nconst body = context.document.body;\n\nconst sentences = body.search(\n  \"Below are the cross-validation results for all the three algorithms.\",\n  { matchCase: true }\n);\nsentences.load(\"items\");\nawait context.sync();\n\nif (sentences.items.length === 0) {\n  throw new Error(\"Target sentence not found.\");\n}\nconst sentence = sentences.items[0];\n\n// Find \"all the three \" inside that sentence only, so the other\n// occurrence of the same phrase elsewhere in the document is left alone.\nconst phrase = sentence.search(\"all the three \", { matchCase: true });\nphrase.load(\"items\");\nawait context.sync();\n\nif (phrase.items.length === 0) {\n  throw new Error(\"Phrase to replace not found inside target sentence.\");\n}\n\n// Replace \"all the three \" with \"both the \" in place; the inserted text\n// inherits the surrounding (Comic Sans MS) run formatting.\nphrase.items[0].insertText(\"both the \", \"Replace\");\nawait context.sync();\n\n// Re-locate the freshly inserted \"both the \" text and nudge its direct\n// character formatting (apply then immediately revert a font toggle).\n// This mirrors how Word itself splits a retyped phrase into its own\n// run(s) even though the resulting formatting is unchanged, matching\n// the surrounding runs exactly.\nconst updated = body.search(\n  \"Below are the cross-validation results for both the algorithms.\",\n  { matchCase: true }\n);\nupdated.load(\"items\");\nawait context.sync();\n\nconst newPhrase = updated.items[0].search(\"both the \", { matchCase: true });\nnewPhrase.load(\"items\");\nawait context.sync();\n\nconst newRange = newPhrase.items[0];\nnewRange.font.bold = true;\nawait context.sync();\nnewRange.font.bold = false;\nawait context.sync();\n", "ps1": "# Change: \"Below are the cross-validation results for all the three\n# algorithms.\" -> \"Below are the cross-validation results for both the\n# algorithms.\" (the sentence right before the algorithm comparison\n# table). The document also contains an unrelated sentence (\"Columns\n# mentioned in all the three categories ...\") that must be left alone,\n# so we locate the specific paragraph that holds the whole target\n# sentence before doing any Find/Replace.\n\n$d = $word.ActiveDocument\n\n$targetIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    if ($para.Range.Text -like \"*cross-validation results for all the three algorithms*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Target sentence not found.\"\n}\n\n$targetRange = $d.Paragraphs.Item($targetIndex).Range\n\n$find = $targetRange.Find\n$find.Text = \"all the three \"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Phrase to replace not found inside target paragraph.\"\n}\n\n# $targetRange now spans just the found \"all the three \" text (Find\n# collapses/extends its owning Range to the match). Replace it in\n# place so it keeps the surrounding (Comic Sans MS) run formatting.\n$targetRange.Text = \"both the \"\n\n# Nudge the direct character formatting of the freshly written text\n# (apply, then immediately revert, a bold toggle). This mirrors how\n# Word splits a retyped phrase into its own run even though the\n# resulting formatting is unchanged, matching the surrounding runs.\n$targetRange.Bold = 1\n$targetRange.Bold = 0\n"}
